$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D4:D5 ("operation" / "[kW_el*h/a]" / "[kW_el*h]") for conversion_1
$ws.Range("D4").Value = 8723.919156434198
$ws.Range("D5").Value = 8723.919156434198

# Update D9:D10 ("operation" / "[kW_el*h/a]" / "[kW_el*h]") for conversion_2
$ws.Range("D9").Value = 11647.48159132677
$ws.Range("D10").Value = 11647.48159132677

# Update D14:D15 ("operation" / "[kW_el*h/a]" / "[kW_el*h]") for conversion_3
$ws.Range("D14").Value = 11493.68084356613
$ws.Range("D15").Value = 11493.68084356613
